$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels (column I/J) on top table ---
$ws.Range("I1").Value = "Av"
$ws.Range("J1").Value = "sigma_Av"

# --- Av / sigma_Av columns for the first measurement table (rows 2-11) ---
$ws.Range("I2").Formula = "=C2/A2"
$ws.Range("I3:I11").Formula = "=C3/A3"

$ws.Range("J2").Formula = "=I2*SQRT((G2/A2)^2 + (H2/C2)^2)"
$ws.Range("J3:J11").Formula = "=I3*SQRT((G3/A3)^2 + (H3/C3)^2)"

# --- New small data block headers (columns M, O, Q) ---
$ws.Range("M1").Value = "#V_in [mV]"
$ws.Range("O1").Value = "V_out [mV]"
$ws.Range("Q1").Value = "Av"

# --- Raw measurement data for the new block (columns M,N,O,P rows 2-13) ---
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 0.4
$ws.Range("O2").Value = 256
$ws.Range("P2").Value = 2

$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 0.8
$ws.Range("O3").Value = 511
$ws.Range("P3").Value = 4

$ws.Range("M4").Value = 150
$ws.Range("N4").Value = 1.2
$ws.Range("O4").Value = 767
$ws.Range("P4").Value = 6

$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 1.6
$ws.Range("O5").Value = 1022
$ws.Range("P5").Value = 8

$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 1278
$ws.Range("P6").Value = 11

$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 1534
$ws.Range("P7").Value = 12

$ws.Range("M8").Value = 349
$ws.Range("N8").Value = 3
$ws.Range("O8").Value = 1790
$ws.Range("P8").Value = 14

$ws.Range("M9").Value = 399
$ws.Range("N9").Value = 3
$ws.Range("O9").Value = 2046
$ws.Range("P9").Value = 16

$ws.Range("M10").Value = 449
$ws.Range("N10").Value = 4
$ws.Range("O10").Value = 2302
$ws.Range("P10").Value = 18

$ws.Range("M11").Value = 499
$ws.Range("N11").Value = 4
$ws.Range("O11").Value = 2558
$ws.Range("P11").Value = 20

$ws.Range("M12").Value = 549
$ws.Range("N12").Value = 4
$ws.Range("O12").Value = 2815
$ws.Range("P12").Value = 22

$ws.Range("M13").Value = 599
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 3073
$ws.Range("P13").Value = 24

# M2:M5 display with one decimal (new custom number format "0.0")
$ws.Range("M2:M5").NumberFormat = "0.0"

# --- Gain / sigma_gain formulas for the new block (columns Q,R rows 2-13) ---
$ws.Range("Q2").Formula = "=O2/M2"
$ws.Range("Q3:Q13").Formula = "=O3/M3"

$ws.Range("R2").Formula = "=Q2*SQRT((N2/M2)^2 + (P2/O2)^2)"
$ws.Range("R3:R13").Formula = "=Q3*SQRT((N3/M3)^2 + (P3/O3)^2)"

# --- View state: scroll back to top, select Q3 ---
$ws.Range("A1").Select()
$ws.Range("Q3").Select()
